$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.302.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.868.34"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.39%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.99"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.44%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4699"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.48%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2865"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06567"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.53"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08023"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.85"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.54%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.870.97"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.111"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6836"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "267.94"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.330.25"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.93"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007624"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.27%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.117.34"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.42%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.257"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.197"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.378"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.49"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.84"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.945"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.79%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09867"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.367"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.460"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.058"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04690"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.132"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6998"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.08%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01866"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.22%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.273"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.96"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.953"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8409"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.86%  "

$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4160"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.83"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.176"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.040"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "906.16"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.39"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05682"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.79%  "
